# Add three new log entries (rows 18-20) to the worklog sheet, matching
# the existing table's layout/styles: date in col A, time-spent in col B,
# description (wrapped) in col C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: 2023-11-13, "~3 hrs", crx.c / TuringCraft meeting note ---
$ws.Range("A18").Value = 45243
$ws.Range("B18").Value = "~3 hrs"
$ws.Range("C18").Value = "cleaned up the code by moving around functions and adding the crx.c and .h files. Also had a 40 min meeting with Professor Weiss about TuringCraft"

# --- Row 19: 2023-11-14, "~1 hr", redesigned plan note ---
$ws.Range("A19").Value = 45244
$ws.Range("B19").Value = "~1 hr"
$ws.Range("C19").Value = "redesigned some of the plan and moved around and commented some more code"

# --- Row 20: 2016-11-15, "~3-4 hrs", vscode extension note ---
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = "~3-4 hrs"
$ws.Range("C20").Value = "looked into and wrote a small extention for vs code to do syntax highlighting for .sul files. Thank you documentation and the ""yo generator-code"" tool"

# Match the description column's wrap/vertical-center style used throughout
# the rest of the table (copy format from an existing wrapped cell).
$ws.Range("C2").Copy()
$ws.Range("C18:C20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# These rows wrap onto several lines in real Excel; reproduce the resulting
# row heights explicitly.
$ws.Rows(18).RowHeight = 75
$ws.Rows(19).RowHeight = 45
$ws.Rows(20).RowHeight = 75

# Restore the frozen header pane and move the live selection to match where
# the author ended up (B21, just past the new last row).
$win = $excel.ActiveWindow
[void]($win.FreezePanes = $false)
$ws.Range("A2").Select() | Out-Null
[void]($win.FreezePanes = $true)
$ws.Range("B21").Select() | Out-Null
